# Add a new "India" data row to the country key table (tables/TableS4_country_key.xlsx).
# India is alphabetically between "China" (row 16) and "Laos" (row 17), so we insert a
# new row at row 17 - which pushes Laos..Sweden down by one row (17-31 -> 18-32) - and
# then populate the new row 17 with India's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing rows 17-31 down to 18-32, leaving a blank row 17 for India.
$ws.Rows.Item(17).Insert() | Out-Null

# Fill in the new row with India's values (Asia / India / IND / Regional / Medium / 0.645 /
# Moderate absorption / Semi-unrefined diet) - matching the pattern of its neighboring rows.
$ws.Cells.Item(17, 1).Value = "Asia"
$ws.Cells.Item(17, 2).Value = "India"
$ws.Cells.Item(17, 3).Value = "IND"
$ws.Cells.Item(17, 4).Value = "Regional"
$ws.Cells.Item(17, 5).Value = "Medium"
$ws.Cells.Item(17, 6).Value = 0.645
$ws.Cells.Item(17, 7).Value = "Moderate absorption"
$ws.Cells.Item(17, 8).Value = "Semi-unrefined diet"

# Restore the view: scroll so row 3 is at the top and select E12 (matches the saved
# workbook view state after the edit).
$excel.ActiveWindow.ScrollRow = 3
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("E12").Select() | Out-Null
